# Apply the 8.9.1.1 tourism GDP table update:
#  - add a new "2022" column (S) to the year header row (row 4)
#  - add the corresponding 2022 data point (S5) to the data row (row 5)
#  - revise the existing 2019-2021 data points (P5, Q5, R5)
#  - move the active selection to T4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create column S by copying the formatting of column R (the previous last
# column) one cell to the right, then overwrite with the new values. Using
# Copy(destination) duplicates the source cell's style along with its value
# so the new cells inherit the same borders/fonts/number formats as the rest
# of the table.
$ws.Range("R4").Copy($ws.Range("S4"))
$ws.Range("S4").Value = 2022

$ws.Range("R5").Copy($ws.Range("S5"))
$ws.Range("S5").Value = 3.4

# Revise the existing values in the data row.
$ws.Range("P5").Value = 4.4
$ws.Range("Q5").Value = 2.9
$ws.Range("R5").Value = 3.2

# Update the current selection to match the author's final cursor position.
$null = $ws.Range("T4").Select()
